$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: B4, C4, D4 -> become text values (shared strings)
$ws.Range("B4").Value = "7558829"
$ws.Range("C4").Value = "75398"
$ws.Range("D4").Value = "3698547"

# Row 7: B7, C7, D7 -> become text values (shared strings)
$ws.Range("B7").Value = "3108228425"
$ws.Range("C7").Value = "3121715639"
$ws.Range("D7").Value = "312321666"

# Update selection to A6
$ws.Range("A6").Select()
